$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

# --- Title placeholder (shape 1): reposition/resize and set the closing text ---
$title = $s.Shapes.Item(1)

$title.Left = 66.0
$title.Top = 217.8125
$title.Width = 828.0
$title.Height = 104.3751

$tr = $title.TextFrame.TextRange
$tr.Text = "Ďakujem"
$tr.LanguageID = "en-GB"
$tr.Font.Bold = $true

$r2 = $tr.InsertAfter(" za ")
$r2.LanguageID = "en-GB"
$r2.Font.Bold = $true

$r3 = $tr.InsertAfter("pozornosť")
$r3.LanguageID = "en-GB"
$r3.Font.Bold = $true

$full = $title.TextFrame.TextRange
$full.ParagraphFormat.Alignment = 2

# --- Content placeholder (shape 2): remove entirely ---
$s.Shapes.Item(2).Delete()
